$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (row 1 uses the bold/bordered/centered header style).
# Set the values first, then copy just the formatting (PasteSpecial formats,
# xlPasteFormats = -4122) from an existing header cell so the underlying
# cell style index is reused rather than a new one being created.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in I2:I18 with 1, and J2:J18 with the same value as the corresponding
# H column cell (mirrors H into J), for each data row.
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 9).Value = 1                                  # column I
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2      # column J
}
